# Refresh the hourly cryptocurrency price / volume snapshot on Sheet1.
# (table rows 2-51, columns B:E), matching the GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Switch the whole table to Text format before writing any values so strings
# such as "1.00", "538.49" or "5.35" are stored as literal text instead of
# being reinterpreted as numbers by Excel's automatic type detection.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '59.291.15'
$ws.Range('E2').Value = '  +0.76%  '
$ws.Range('D3').Value = '2.526.91'
$ws.Range('E3').Value = '  +0.59%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '538.49'
$ws.Range('E5').Value = '  +1.49%  '
$ws.Range('D6').Value = '138.08'
$ws.Range('E6').Value = '  -0.28%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '0.568'
$ws.Range('E8').Value = '  +1.18%  '
$ws.Range('D9').Value = '2.526.96'
$ws.Range('E9').Value = '  +0.55%  '
$ws.Range('E10').Value = '  +1.56%  '
$ws.Range('E11').Value = '  -0.94%  '
$ws.Range('D12').Value = '5.35'
$ws.Range('E12').Value = '  -1.81%  '
$ws.Range('E13').Value = '  -2.08%  '
$ws.Range('D14').Value = '2.974.76'
$ws.Range('E14').Value = '  +0.63%  '
$ws.Range('D15').Value = '23.19'
$ws.Range('E15').Value = '  +0.56%  '
$ws.Range('D16').Value = '59.229.55'
$ws.Range('E16').Value = '  +0.73%  '
$ws.Range('D18').Value = '2.521.08'
$ws.Range('E18').Value = '  +0.53%  '
$ws.Range('D19').Value = '11.14'
$ws.Range('E19').Value = '  +1.09%  '
$ws.Range('E20').Value = '  +0.69%  '
$ws.Range('D21').Value = '325.73'
$ws.Range('E21').Value = '  +1.05%  '
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('E23').Value = '  +2.94%  '
$ws.Range('D24').Value = '65.56'
$ws.Range('E24').Value = '  +5.32%  '
$ws.Range('E25').Value = '  -0.27%  '
$ws.Range('E26').Value = '  +0.73%  '
$ws.Range('E27').Value = '  +0.37%  '
$ws.Range('D28').Value = '7.68'
$ws.Range('E28').Value = '  -1.03%  '
$ws.Range('D29').Value = '0.0₃0779'
$ws.Range('E29').Value = '  +1.31%  '
$ws.Range('D30').Value = '6.71'
$ws.Range('E30').Value = '  +0.55%  '
$ws.Range('E31').Value = '  +0.26%  '
$ws.Range('E32').Value = '  +6.45%  '
$ws.Range('D33').Value = '165.75'
$ws.Range('E33').Value = '  +1.59%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = '1.47'
$ws.Range('E34').Value = '  +3.73%  '
$ws.Range('B35').Value = 'USDe'
$ws.Range('C35').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D35').Value = '0.997'
$ws.Range('E35').Value = '  -0.15%  '
$ws.Range('D36').Value = '18.50'
$ws.Range('E36').Value = '  +0.37%  '
$ws.Range('E37').Value = '  -2.08%  '
$ws.Range('E38').Value = '  -0.33%  '
$ws.Range('E39').Value = '  +0.15%  '
$ws.Range('D40').Value = '0.821'
$ws.Range('E40').Value = '  +2.44%  '
$ws.Range('D41').Value = '3.65'
$ws.Range('E41').Value = '  +0.14%  '
$ws.Range('D42').Value = '290.23'
$ws.Range('E42').Value = '  +4.11%  '
$ws.Range('D43').Value = '5.24'
$ws.Range('E43').Value = '  +0.38%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  +0.13%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = '131.72'
$ws.Range('E45').Value = '  +8.38%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').Value = '0.611'
$ws.Range('E46').Value = '  +2.51%  '
$ws.Range('E47').Value = '  +0.29%  '
$ws.Range('D49').Value = '0.0511'
$ws.Range('E49').Value = '  +0.37%  '
$ws.Range('E50').Value = '  -0.31%  '
$ws.Range('D51').Value = '17.46'
$ws.Range('E51').Value = '  -0.96%  '

# Restore the default (unstyled) cell style now that the literal text values
# are locked in, so the edited cells look exactly like their neighbours.
$dataRange.Style = "Normal"
